$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

# New booster translation rows to append after the existing "lab.base.*" block
# (row 354 is the last existing data row, new rows start at 355)
$rows = @(
    @("lab.booster.tooltip.create", "Vytvořit booster"),
    @("lab.booster.create.title", "Nový booster"),
    @("lab.booster.create.subtitle", "Boostery jsou užitečné pro tvorbu mixů."),
    @("lab.booster.name.label", "Název"),
    @("lab.booster.vendorId.label", "Výrobce"),
    @("lab.booster.nicotine.label", "Obsah nikotinu"),
    @("lab.booster.volume.label", "Objem"),
    @("lab.booster.create.submit", "Vytvořit booster"),
    @("lab.booster.create.success", "Booster [{{data.name}}] vytvořen."),
    @("lab.booster.pg.label", "PG"),
    @("lab.booster.vg.label", "VG")
)

$startRow = 355
$templateRow = 354

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = "cs"
    $ws.Cells.Item($r, 2).Value = $rows[$i][0]
    $ws.Cells.Item($r, 3).Value = $rows[$i][1]

    # Copy the formatting of the last existing data row onto the new row
    # so the new cells keep the same style index ("import" cell style)
    # instead of falling back to the column's default style.
    $ws.Range("A$templateRow`:C$templateRow").Copy() | Out-Null
    $ws.Range("A$r`:C$r").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false

$lastRow = $startRow + $rows.Length - 1

# Match the new selection recorded in the saved workbook
$ws.Range("B360").Select()
